$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells (changed F/G values across many rows)
$ws.Cells.Item(313,6).Value = 76823
$ws.Cells.Item(321,6).Value = 91100
$ws.Cells.Item(322,6).Value = 110377
$ws.Cells.Item(323,6).Value = 217375
$ws.Cells.Item(327,6).Value = 225530
$ws.Cells.Item(328,6).Value = 181015
$ws.Cells.Item(330,6).Value = 72253
$ws.Cells.Item(331,6).Value = 154367
$ws.Cells.Item(337,6).Value = 105615
$ws.Cells.Item(338,6).Value = 221348
$ws.Cells.Item(341,6).Value = 286154
$ws.Cells.Item(342,6).Value = 178918
$ws.Cells.Item(343,6).Value = 133993
$ws.Cells.Item(344,6).Value = 136295
$ws.Cells.Item(345,6).Value = 293836
$ws.Cells.Item(345,7).Value = 3343
$ws.Cells.Item(348,6).Value = 234188
$ws.Cells.Item(349,6).Value = 161088
$ws.Cells.Item(350,6).Value = 128036
$ws.Cells.Item(350,7).Value = 2796
$ws.Cells.Item(351,6).Value = 150262
$ws.Cells.Item(352,6).Value = 307955
$ws.Cells.Item(352,7).Value = 3575
$ws.Cells.Item(354,6).Value = 317013
$ws.Cells.Item(355,6).Value = 222343
$ws.Cells.Item(356,6).Value = 160796
$ws.Cells.Item(357,6).Value = 138320
$ws.Cells.Item(358,6).Value = 159469
$ws.Cells.Item(359,6).Value = 321000
$ws.Cells.Item(359,7).Value = 3334
$ws.Cells.Item(362,6).Value = 228978
$ws.Cells.Item(363,6).Value = 189551
$ws.Cells.Item(364,6).Value = 168765
$ws.Cells.Item(365,6).Value = 185186
$ws.Cells.Item(365,7).Value = 2400
$ws.Cells.Item(366,6).Value = 339253
$ws.Cells.Item(366,7).Value = 2841
$ws.Cells.Item(367,6).Value = 769233
$ws.Cells.Item(367,7).Value = 3921
$ws.Cells.Item(368,6).Value = 346318
$ws.Cells.Item(368,7).Value = 2298
$ws.Cells.Item(369,6).Value = 235361
$ws.Cells.Item(370,6).Value = 180672
$ws.Cells.Item(371,6).Value = 160315
$ws.Cells.Item(372,6).Value = 178700
$ws.Cells.Item(373,6).Value = 350753
$ws.Cells.Item(374,6).Value = 774878
$ws.Cells.Item(374,7).Value = 3431
$ws.Cells.Item(375,6).Value = 350345
$ws.Cells.Item(376,6).Value = 222838
$ws.Cells.Item(376,7).Value = 2230
$ws.Cells.Item(377,6).Value = 176858
$ws.Cells.Item(378,6).Value = 157508
$ws.Cells.Item(379,6).Value = 180839
$ws.Cells.Item(380,6).Value = 345746
$ws.Cells.Item(380,7).Value = 2035
$ws.Cells.Item(381,6).Value = 748087
$ws.Cells.Item(381,7).Value = 2691
$ws.Cells.Item(382,6).Value = 357396
$ws.Cells.Item(383,6).Value = 222696
$ws.Cells.Item(384,6).Value = 172356
$ws.Cells.Item(385,6).Value = 151086
$ws.Cells.Item(386,6).Value = 183268
$ws.Cells.Item(387,6).Value = 351679
$ws.Cells.Item(388,6).Value = 731148
$ws.Cells.Item(388,7).Value = 2203
$ws.Cells.Item(390,6).Value = 220020
$ws.Cells.Item(391,6).Value = 178071
$ws.Cells.Item(392,6).Value = 222028
$ws.Cells.Item(392,7).Value = 1240
$ws.Cells.Item(393,6).Value = 309294
$ws.Cells.Item(395,6).Value = 753427
$ws.Cells.Item(398,6).Value = 300515
$ws.Cells.Item(400,6).Value = 149874
$ws.Cells.Item(400,7).Value = 805
$ws.Cells.Item(401,6).Value = 273034
$ws.Cells.Item(404,6).Value = 224883
$ws.Cells.Item(408,6).Value = 305856
$ws.Cells.Item(408,7).Value = 836
$ws.Cells.Item(410,6).Value = 365106
$ws.Cells.Item(411,6).Value = 225445
$ws.Cells.Item(412,6).Value = 176886
$ws.Cells.Item(414,6).Value = 149159
$ws.Cells.Item(415,6).Value = 308816
$ws.Cells.Item(416,6).Value = 674083
$ws.Cells.Item(416,7).Value = 934
$ws.Cells.Item(417,6).Value = 344480
$ws.Cells.Item(417,7).Value = 590
$ws.Cells.Item(418,6).Value = 202513
$ws.Cells.Item(419,6).Value = 149713
$ws.Cells.Item(420,6).Value = 139100
$ws.Cells.Item(421,6).Value = 153357
$ws.Cells.Item(423,6).Value = 440510
$ws.Cells.Item(423,7).Value = 638
$ws.Cells.Item(424,6).Value = 266848
$ws.Cells.Item(425,6).Value = 138171
$ws.Cells.Item(426,6).Value = 107417
$ws.Cells.Item(428,6).Value = 102554
$ws.Cells.Item(429,6).Value = 178476
$ws.Cells.Item(430,6).Value = 175556
$ws.Cells.Item(436,6).Value = 145433
$ws.Cells.Item(440,6).Value = 73746
$ws.Cells.Item(443,6).Value = 106931
$ws.Cells.Item(446,6).Value = 86699
$ws.Cells.Item(447,6).Value = 67038
$ws.Cells.Item(447,7).Value = 207
$ws.Cells.Item(449,6).Value = 60011
$ws.Cells.Item(449,7).Value = 157
$ws.Cells.Item(450,6).Value = 91767
$ws.Cells.Item(456,6).Value = 50373
$ws.Cells.Item(457,6).Value = 79003
$ws.Cells.Item(458,6).Value = 70696
$ws.Cells.Item(459,6).Value = 59829
$ws.Cells.Item(460,6).Value = 58421
$ws.Cells.Item(461,6).Value = 45323
$ws.Cells.Item(462,6).Value = 43616
$ws.Cells.Item(462,7).Value = 50
$ws.Cells.Item(464,6).Value = 73612
$ws.Cells.Item(465,6).Value = 61758
$ws.Cells.Item(465,7).Value = 58
$ws.Cells.Item(466,6).Value = 51148
$ws.Cells.Item(467,6).Value = 52272
$ws.Cells.Item(468,6).Value = 41828
$ws.Cells.Item(469,6).Value = 41121
$ws.Cells.Item(470,6).Value = 43521
$ws.Cells.Item(471,6).Value = 66821
$ws.Cells.Item(472,6).Value = 51769
$ws.Cells.Item(473,6).Value = 39974
$ws.Cells.Item(473,7).Value = 41
$ws.Cells.Item(474,6).Value = 45297
$ws.Cells.Item(475,6).Value = 36565
$ws.Cells.Item(476,6).Value = 37223
$ws.Cells.Item(476,7).Value = 31
$ws.Cells.Item(477,6).Value = 37025
$ws.Cells.Item(477,7).Value = 37
$ws.Cells.Item(478,6).Value = 54686
$ws.Cells.Item(479,6).Value = 42318
$ws.Cells.Item(480,6).Value = 33362
$ws.Cells.Item(481,6).Value = 41365
$ws.Cells.Item(482,6).Value = 35903
$ws.Cells.Item(483,6).Value = 64657
$ws.Cells.Item(488,6).Value = 6072
$ws.Cells.Item(489,6).Value = 12563
$ws.Cells.Item(490,6).Value = 10681
$ws.Cells.Item(491,6).Value = 9877
$ws.Cells.Item(492,6).Value = 13830
$ws.Cells.Item(493,6).Value = 8126
$ws.Cells.Item(494,6).Value = 6278
$ws.Cells.Item(495,6).Value = 10219
$ws.Cells.Item(496,6).Value = 8016
$ws.Cells.Item(497,6).Value = 7561
$ws.Cells.Item(498,6).Value = 9023
$ws.Cells.Item(499,6).Value = 10662
$ws.Cells.Item(499,7).Value = 11
$ws.Cells.Item(500,6).Value = 7079
$ws.Cells.Item(500,7).Value = 8
$ws.Cells.Item(501,6).Value = 5615
$ws.Cells.Item(501,7).Value = 7

# Add new rows 502-504
$ws.Cells.Item(502,1).Value = 44396
$ws.Cells.Item(502,1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(502,2).Value = 392139
$ws.Cells.Item(502,3).Value = 7892
$ws.Cells.Item(502,4).Value = 35
$ws.Cells.Item(502,5).Value = 12531
$ws.Cells.Item(502,6).Value = 9379
$ws.Cells.Item(502,7).Value = 17

$ws.Cells.Item(503,1).Value = 44397
$ws.Cells.Item(503,1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(503,2).Value = 392185
$ws.Cells.Item(503,3).Value = 6603
$ws.Cells.Item(503,4).Value = 46
$ws.Cells.Item(503,5).Value = 12534
$ws.Cells.Item(503,6).Value = 6832
$ws.Cells.Item(503,7).Value = 6

$ws.Cells.Item(504,1).Value = 44398
$ws.Cells.Item(504,1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(504,2).Value = 392219
$ws.Cells.Item(504,3).Value = 6779
$ws.Cells.Item(504,4).Value = 34
$ws.Cells.Item(504,5).Value = 12534
$ws.Cells.Item(504,6).Value = 5385
$ws.Cells.Item(504,7).Value = 15
